$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> (DAMSLTag, DialogAct) following SGNN re-annotation
$updates = @(
    @{ Row = 6; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 25; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 28; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 29; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 32; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 38; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 64; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 65; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 76; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 77; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 81; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 82; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 95; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 97; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 100; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 103; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 104; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 124; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 147; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 150; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 173; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 178; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 180; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 181; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 184; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 187; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 189; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 190; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 198; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 203; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 206; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 207; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 219; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 220; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 229; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 236; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 238; DAMSLTag = '%'; DialogAct = 'Uninterpretable' }
    @{ Row = 266; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 271; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 279; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 290; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 291; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 317; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 324; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 325; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 330; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 348; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 349; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 352; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 353; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 356; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.DAMSLTag
    $ws.Cells.Item($u.Row, 10).Value = $u.DialogAct
}

$wb.Save()
